$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 3 values: C3 and D3 now point to "hello" (same string as B3), E3 becomes 2
$ws.Range("C3").Value = "hello"
$ws.Range("D3").Value = "hello"
$ws.Range("E3").Value = 2

# Remove row 4 entirely (shifts nothing below it up, it's the last row)
$ws.Rows.Item(4).Delete()
